$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("V4").Value = '750 (S2022T03) + 750 (S2022T01)'
$ws.Range("V7").Value = '750 (S2022T02) + 500 (S2022T03) + 500 (S2022T01)'
$ws.Range("F21").Value = 5
$ws.Range("Z21").Value = 'Dupertuis, Gaston'
$ws.Range("F22").Value = 11
$ws.Range("Z22").Value = 'Nowotny, Martin'
$ws.Range("V25").Value = '500 (S2022T02) + 500 (S2022T01) + 250 (S2022T03)'
$ws.Range("F38").Value = 15
$ws.Range("Z38").Value = 'Rodriguez, Pablo'
$ws.Range("F39").Value = 31
$ws.Range("Z39").Value = 'Fernandez, Carlos'
$ws.Range("F62").Value = 94
$ws.Range("Z62").Value = 'Sartor, Yemel'
$ws.Range("F63").Value = 145
$ws.Range("Z63").Value = 'Sueldo, Pablo'
$ws.Range("F80").Value = 158
$ws.Range("Z80").Value = 'Campos, Dario'
$ws.Range("F81").Value = 160
$ws.Range("Z81").Value = 'Chiara, Lucio'
$ws.Range("F91").Value = 16
$ws.Range("Z91").Value = 'Rulfi, Daniel'
$ws.Range("F92").Value = 196
$ws.Range("Z92").Value = 'Escobar, Esteban'
$ws.Range("F95").Value = 57
$ws.Range("Z95").Value = 'Apaza, Catriel'
$ws.Range("F96").Value = 122
$ws.Range("Z96").Value = 'Botta, Maximiliano'
$ws.Range("F100").Value = 59
$ws.Range("Z100").Value = 'Arrieta, Maximiliano'
$ws.Range("F101").Value = 124
$ws.Range("Z101").Value = 'Mendieta, Elias'
$ws.Range("F102").Value = 127
$ws.Range("Z102").Value = 'Badano, Pablo'
$ws.Range("F105").Value = 194
$ws.Range("Z105").Value = 'Asenie, Santiago'
$ws.Range("F122").Value = 34
$ws.Range("I122").Value = $true
$ws.Range("K122").Value = 125
$ws.Range("O122").Value = 625
$ws.Range("S122").Value = 3
$ws.Range("W122").Value = '250 (S2022T02) + 250 (S2022T01) + 125 (S2022T03)'
$ws.Range("Z122").Value = 'Levin, Raul'
$ws.Range("F124").Value = 223
$ws.Range("I124").Value = $false
$ws.Range("K124").Value = 0
$ws.Range("O124").Value = 0
$ws.Range("S124").Value = 0
$ws.Range("W124").Value = ""
$ws.Range("Z124").Value = 'Aguer, Jose'
$ws.Range("X127").Value = '125 (S2022T03) + 125 (S2022T01) + 65 (S2022T02)'
$ws.Range("F137").Value = 121
$ws.Range("Z137").Value = 'Tenca, Javier'
$ws.Range("F138").Value = 211
$ws.Range("Z138").Value = 'Miner, Alberto'
$ws.Range("F150").Value = 230
$ws.Range("Z150").Value = 'Godoy, Franco'
$ws.Range("F151").Value = 244
$ws.Range("Z151").Value = 'Fucks, Alyssa'
$ws.Range("F156").Value = 167
$ws.Range("Z156").Value = 'Cossi, Francisco'
$ws.Range("F157").Value = 203
$ws.Range("Z157").Value = 'Brian, Martin'
$ws.Range("F187").Value = 72
$ws.Range("K187").Value = 75
$ws.Range("O187").Value = 75
$ws.Range("P187").Value = 295
$ws.Range("S187").Value = 1
$ws.Range("W187").Value = '75 (S2022T03)'
$ws.Range("X187").Value = '190 (S2022T01) + 65 (S2022T03) + 40 (S2022T02)'
$ws.Range("Z187").Value = 'Colavini, Daniel'
$ws.Range("F188").Value = 280
$ws.Range("K188").Value = 0
$ws.Range("O188").Value = 0
$ws.Range("P188").Value = 195
$ws.Range("S188").Value = 0
$ws.Range("W188").Value = ""
$ws.Range("X188").Value = '65 (S2022T03) + 65 (S2022T02) + 65 (S2022T01)'
$ws.Range("Z188").Value = 'Savino, Leandro'
$ws.Range("F197").Value = 204
$ws.Range("Z197").Value = 'Delgado, Pablo'
$ws.Range("F198").Value = 240
$ws.Range("Z198").Value = 'Arrieta, Matias'
$ws.Range("F199").Value = 249
$ws.Range("Z199").Value = 'Muller, Tomas'
$ws.Range("F200").Value = 256
$ws.Range("Z200").Value = 'Portillo, Lucas'
$ws.Range("F205").Value = 213
$ws.Range("Z205").Value = 'Pillac, Juan Pablo'
$ws.Range("F206").Value = 266
$ws.Range("Z206").Value = 'Comas, Javier'
$ws.Range("F215").Value = 197
$ws.Range("I215").Value = $true
$ws.Range("L215").Value = 40
$ws.Range("O215").Value = 50
$ws.Range("P215").Value = 230
$ws.Range("S215").Value = 1
$ws.Range("T215").Value = 3
$ws.Range("W215").Value = '50 (S2022T01)'
$ws.Range("X215").Value = '125 (S2022T02) + 65 (S2022T01) + 40 (S2022T03)'
$ws.Range("Z215").Value = 'Larrosa, Jorge'
$ws.Range("F216").Value = 229
$ws.Range("I216").Value = $false
$ws.Range("L216").Value = 0
$ws.Range("O216").Value = 0
$ws.Range("P216").Value = 0
$ws.Range("S216").Value = 0
$ws.Range("T216").Value = 0
$ws.Range("W216").Value = ""
$ws.Range("X216").Value = ""
$ws.Range("Z216").Value = 'Aguirre, Sandra'
$ws.Range("F242").Value = 295
$ws.Range("Z242").Value = 'Antunez, Pablo'
$ws.Range("F243").Value = 299
$ws.Range("Z243").Value = 'Ferrero, Alejandro'
$ws.Range("F244").Value = 285
$ws.Range("Z244").Value = 'Lell, Claudia'
$ws.Range("F245").Value = 304
$ws.Range("Z245").Value = 'Velazquez, Pedro'
$ws.Range("F252").Value = 306
$ws.Range("Z252").Value = 'Bertoli, Julian'
$ws.Range("F253").Value = 307
$ws.Range("Z253").Value = 'Bertoli, Maximiliano'
